$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '69.839.88'
$ws.Range('E2').Value = '  +0.27%  '
Set-TextValue $ws.Range('D3') '3.561.67'
$ws.Range('E3').Value = '  -0.47%  '
$ws.Range('E4').Value = '  +0.10%  '
Set-TextValue $ws.Range('D5') '576.95'
$ws.Range('E5').Value = '  -2.51%  '
Set-TextValue $ws.Range('D6') '188.20'
$ws.Range('E6').Value = '  -1.64%  '
Set-TextValue $ws.Range('D7') '0.630'
$ws.Range('E7').Value = '  -2.72%  '
Set-TextValue $ws.Range('D8') '3.559.67'
$ws.Range('E8').Value = '  -0.37%  '
$ws.Range('E9').Value = '  +0.01%  '
$ws.Range('E10').Value = '  -2.66%  '
Set-TextValue $ws.Range('D11') '0.657'
$ws.Range('E11').Value = '  -0.48%  '
Set-TextValue $ws.Range('D12') '56.01'
$ws.Range('E12').Value = '  -3.33%  '
$ws.Range('E13').Value = '  +1.82%  '
Set-TextValue $ws.Range('D14') '9.67'
$ws.Range('E14').Value = '  -0.43%  '
Set-TextValue $ws.Range('D15') '4.135.53'
$ws.Range('E15').Value = '  -0.52%  '
Set-TextValue $ws.Range('D16') '19.89'
$ws.Range('E16').Value = '  +2.83%  '
Set-TextValue $ws.Range('D17') '3.565.51'
$ws.Range('E17').Value = '  -0.32%  '
Set-TextValue $ws.Range('D18') '69.734.58'
$ws.Range('E18').Value = '  +0.18%  '
Set-TextValue $ws.Range('D19') '12.48'
$ws.Range('E19').Value = '  -1.28%  '
$ws.Range('E20').Value = '  +0.39%  '
Set-TextValue $ws.Range('D21') '1.04'
$ws.Range('E21').Value = '  -1.28%  '
Set-TextValue $ws.Range('D22') '471.64'
$ws.Range('E22').Value = '  -5.78%  '
Set-TextValue $ws.Range('D23') '19.35'
$ws.Range('E23').Value = '  +13.88%  '
Set-TextValue $ws.Range('D24') '5.07'
$ws.Range('E24').Value = '  -7.17%  '
Set-TextValue $ws.Range('D25') '4.34'
$ws.Range('E25').Value = '  -2.63%  '
Set-TextValue $ws.Range('D26') '88.29'
$ws.Range('E26').Value = '  -3.00%  '
Set-TextValue $ws.Range('D27') '3.05'
$ws.Range('E27').Value = '  -1.02%  '
Set-TextValue $ws.Range('D28') '10.94'
$ws.Range('E28').Value = '  -1.51%  '
Set-TextValue $ws.Range('D29') '9.36'
$ws.Range('E29').Value = '  +0.20%  '
Set-TextValue $ws.Range('D30') '32.01'
$ws.Range('E30').Value = '  -0.50%  '
Set-TextValue $ws.Range('D31') '7.63'
$ws.Range('E31').Value = '  +1.96%  '
$ws.Range('E32').Value = '  +2.87%  '
Set-TextValue $ws.Range('D33') '12.07'
$ws.Range('E33').Value = '  -0.79%  '
Set-TextValue $ws.Range('D34') '65.55'
$ws.Range('E34').Value = '  +0.33%  '
Set-TextValue $ws.Range('D35') '583.49'
$ws.Range('E35').Value = '  -5.15%  '
Set-TextValue $ws.Range('D36') '38.74'
$ws.Range('E36').Value = '  +2.02%  '
Set-TextValue $ws.Range('D37') '1.00'
$ws.Range('E37').Value = '  -0.03%  '
$ws.Range('E38').Value = '  -3.38%  '
Set-TextValue $ws.Range('D39') '0.396'
$ws.Range('E39').Value = '  -0.78%  '
$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range('D40') '0.140'
$ws.Range('E40').Value = '  -4.64%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range('D41') '3.50'
$ws.Range('E41').Value = '  -3.51%  '
Set-TextValue $ws.Range('D42') '3.17'
$ws.Range('E42').Value = '  +14.46%  '
$ws.Range('E43').Value = '  +6.63%  '
Set-TextValue $ws.Range('D44') '3.218.36'
$ws.Range('E44').Value = '  -3.80%  '
Set-TextValue $ws.Range('D45') '3.12'
$ws.Range('E45').Value = '  +0.47%  '
Set-TextValue $ws.Range('D46') '0.0443'
$ws.Range('E46').Value = '  -0.20%  '
Set-TextValue $ws.Range('D47') '9.42'
$ws.Range('E47').Value = '  +4.01%  '
$ws.Range('E48').Value = '  +1.54%  '
Set-TextValue $ws.Range('D49') '0.137'
$ws.Range('E49').Value = '  -0.56%  '
Set-TextValue $ws.Range('D50') '0.999'
$ws.Range('E50').Value = '  -0.28%  '
$ws.Range('B51').Value = 'Monero'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range('D51') '137.68'
$ws.Range('E51').Value = '  -2.52%  '
